$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated NATMI edge metrics per "Natmi following Dr Hou advice"
# Each entry: row, then values for columns E,G,H,I,J,K,M,N,O,P,Q,R,S,T
$rowsData = @(
    @(2, 3, 1.842080666666667, 5.526242, 0.1262299803130056, 0.1502126669235156, 3, 44.13164066666667, 132.394922, 0.4415399811720331, 0.4562856844211927, 81.29404206034712, 731.646378543124, 0.0557355831307506, 0.06853988953592899),
    @(3, 3, 1.842080666666667, 5.526242, 0.1262299803130056, 0.1502126669235156, 3, 14.93259333333333, 44.79778, 0.1494015830739255, 0.1543910098595022, 27.50704148252889, 247.56337334276, 0.01885895889015348, 0.02319148534001062),
    @(4, 3, 1.842080666666667, 5.526242, 0.1262299803130056, 0.1502126669235156, 3, 13.61024133333333, 40.830724, 0.1361713639304118, 0.1407189532975654, 25.07116242880089, 225.640461859208, 0.01718890858813099, 0.02113776926151295),
    @(5, 3, 1.842080666666667, 5.526242, 0.1262299803130056, 0.1502126669235156, 3, 17.584752, 52.754256, 0.1759366057935712, 0.1818121982434553, 32.392531687328, 291.532785185952, 0.02220847428565952, 0.02731049517737635),
    @(6, 3, 1.842080666666667, 5.526242, 0.1262299803130056, 0.1502126669235156, 2, 9.6901375, 19.380275, 0.09695046603005844, 0.06679215417828435, 17.85001494609167, 107.10008967655, 0.01223805541831099, 0.01003302760868673),
    @(7, 3, 4.241638, 12.724914, 0.2906614736930972, 0.3458848288425265, 3, 44.13164066666667, 132.394922, 0.4415399811720331, 0.4562856844211927, 187.1904440540787, 1684.713996486708, 0.1283386616218856, 0.1578222958593193),
    @(8, 3, 4.241638, 12.724914, 0.2906614736930972, 0.3458848288425265, 3, 14.93259333333333, 44.79778, 0.1494015830739255, 0.1543910098595022, 63.33865532121334, 570.04789789092, 0.04342528430834888, 0.05340150802007872),
    @(9, 3, 4.241638, 12.724914, 0.2906614736930972, 0.3458848288425265, 3, 13.61024133333333, 40.830724, 0.1361713639304118, 0.1407189532975654, 57.72971682863733, 519.567451457736, 0.03957976931481255, 0.04867255107622789),
    @(10, 3, 4.241638, 12.724914, 0.2906614736930972, 0.3458848288425265, 3, 17.584752, 52.754256, 0.1759366057935712, 0.1818121982434553, 74.588152303776, 671.293370733984, 0.05113799311652092, 0.06288608107092104),
    @(11, 3, 4.241638, 12.724914, 0.2906614736930972, 0.3458848288425265, 2, 9.6901375, 19.380275, 0.09695046603005844, 0.06679215417828435, 41.10205544522501, 246.61233267135, 0.02817976533152935, 0.02310239281597952),
    @(12, 3, 0.6198420000000001, 1.859526, 0.04247514502106894, 0.05054508283814161, 3, 44.13164066666667, 132.394922, 0.4415399811720331, 0.4562856844211927, 27.35464441410801, 246.1917997269721, 0.01875447473288215, 0.02306299771692732),
    @(13, 3, 0.6198420000000001, 1.859526, 0.04247514502106894, 0.05054508283814161, 3, 14.93259333333333, 44.79778, 0.1494015830739255, 0.1543910098595022, 9.255848516920002, 83.30263665228001, 0.006345853907442263, 0.007803706382812875),
    @(14, 3, 0.6198420000000001, 1.859526, 0.04247514502106894, 0.05054508283814161, 3, 13.61024133333333, 40.830724, 0.1361713639304118, 0.1407189532975654, 8.436199208536001, 75.925792876824, 0.005783898430660996, 0.007112651151322025),
    @(15, 3, 0.6198420000000001, 1.859526, 0.04247514502106894, 0.05054508283814161, 3, 17.584752, 52.754256, 0.1759366057935712, 0.1818121982434553, 10.899767849184, 98.097910642656, 0.007472932845596575, 0.009189712621200075),
    @(16, 3, 0.6198420000000001, 1.859526, 0.04247514502106894, 0.05054508283814161, 2, 9.6901375, 19.380275, 0.09695046603005844, 0.06679215417828435, 6.006354208275002, 36.03812524965001, 0.00411798510448695, 0.003376014965879309),
    @(17, 3, 0.8997893333333332, 2.699368, 0.06165874920019015, 0.07337341837147134, 3, 44.13164066666667, 132.394922, 0.4415399811720331, 0.4562856844211927, 39.70917953436622, 357.382615809296, 0.02722480296094307, 0.03347924041994931),
    @(18, 3, 0.8997893333333332, 2.699368, 0.06165874920019015, 0.07337341837147134, 3, 14.93259333333333, 44.79778, 0.1494015830739255, 0.1543910098595022, 13.43618820033778, 120.92569380304, 0.009211914740866545, 0.01132819615921521),
    @(19, 3, 0.8997893333333332, 2.699368, 0.06165874920019015, 0.07337341837147134, 3, 13.61024133333333, 40.830724, 0.1361713639304118, 0.1407189532975654, 12.24634997582578, 110.217149782432, 0.008396155976833079, 0.0103250306330978),
    @(20, 3, 0.8997893333333332, 2.699368, 0.06165874920019015, 0.07337341837147134, 3, 17.584752, 52.754256, 0.1759366057935712, 0.1818121982434553, 15.822572278912, 142.403150510208, 0.01084803105175853, 0.01334018248675393),
    @(21, 3, 0.8997893333333332, 2.699368, 0.06165874920019015, 0.07337341837147134, 2, 9.6901375, 19.380275, 0.09695046603005844, 0.06679215417828435, 8.719082361033333, 52.3144941662, 0.005977844469788928, 0.004900768672455075),
    @(22, 2, 6.989701999999999, 13.979404, 0.478974651772638, 0.3799840030243449, 3, 44.13164066666667, 132.394922, 0.4415399811720331, 0.4562856844211927, 308.4670170310814, 1850.802102186488, 0.2114864587255717, 0.1733812608890678),
    @(23, 2, 6.989701999999999, 13.979404, 0.478974651772638, 0.3799840030243449, 3, 14.93259333333333, 44.79778, 0.1494015830739255, 0.1543910098595022, 104.3743774871867, 626.24626492312, 0.07155957122711432, 0.05866611395738474),
    @(24, 2, 6.989701999999999, 13.979404, 0.478974651772638, 0.3799840030243449, 3, 13.61024133333333, 40.830724, 0.1361713639304118, 0.1407189532975654, 95.13153106808265, 570.7891864084959, 0.06522263161997414, 0.05347095117540476),
    @(25, 2, 6.989701999999999, 13.979404, 0.478974651772638, 0.3799840030243449, 3, 17.584752, 52.754256, 0.1759366057935712, 0.1818121982434553, 122.912176223904, 737.4730573434239, 0.08426917449403568, 0.06908572688720394),
    @(26, 2, 6.989701999999999, 13.979404, 0.478974651772638, 0.3799840030243449, 2, 9.6901375, 19.380275, 0.09695046603005844, 0.06679215417828435, 67.73117346402499, 270.9246938561, 0.04643681570594221, 0.02537995011528371)
)

$colNames = @("E","G","H","I","J","K","M","N","O","P","Q","R","S","T")

foreach ($rowEntry in $rowsData) {
    $rowNum = $rowEntry[0]
    for ($i = 0; $i -lt $colNames.Length; $i++) {
        $col = $colNames[$i]
        $val = $rowEntry[$i + 1]
        $ws.Range("$col$rowNum").Value = $val
    }
}
